# Auto-generated edit script
# 1) Fix comma-separated name fields to use periods (per commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E53").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E99").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E122").Value = "FERNANDEZ. MARIO HUGO"
$ws.Range("E203").Value = "DODERA. JORGE ABELARDO"
$ws.Range("E209").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E222").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# 2) Normalize "Importe" (column H) number-formatted text:
#    remove thousands separators ("."), convert decimal comma to dot.
#    Values are written with a leading apostrophe so Excel keeps them
#    as literal text (matching the original shared-string storage)
#    instead of auto-converting the numeric-looking text into a number.
$ws.Range("H2").Value = "'4660.00"
$ws.Range("H3").Value = "'1040.00"
$ws.Range("H4").Value = "'17640.00"
$ws.Range("H5").Value = "'4155.00"
$ws.Range("H6").Value = "'1340.00"
$ws.Range("H7").Value = "'66000.00"
$ws.Range("H8").Value = "'66000.00"
$ws.Range("H9").Value = "'300000.00"
$ws.Range("H10").Value = "'380000.00"
$ws.Range("H11").Value = "'1154962.09"
$ws.Range("H12").Value = "'6300.00"
$ws.Range("H13").Value = "'994.98"
$ws.Range("H14").Value = "'8000.20"
$ws.Range("H15").Value = "'328.70"
$ws.Range("H16").Value = "'3300.00"
$ws.Range("H17").Value = "'200.00"
$ws.Range("H18").Value = "'87150.00"
$ws.Range("H19").Value = "'98932.40"
$ws.Range("H20").Value = "'14700.00"
$ws.Range("H21").Value = "'127789.90"
$ws.Range("H22").Value = "'7500.00"
$ws.Range("H23").Value = "'102820.90"
$ws.Range("H24").Value = "'49222.00"
$ws.Range("H25").Value = "'58035.81"
$ws.Range("H26").Value = "'3690.00"
$ws.Range("H27").Value = "'3602.90"
$ws.Range("H28").Value = "'40712.59"
$ws.Range("H29").Value = "'17904.00"
$ws.Range("H30").Value = "'24971.27"
$ws.Range("H31").Value = "'28138.00"
$ws.Range("H32").Value = "'9000.00"
$ws.Range("H33").Value = "'2800.00"
$ws.Range("H34").Value = "'8000.00"
$ws.Range("H35").Value = "'7790.00"
$ws.Range("H36").Value = "'525.09"
$ws.Range("H37").Value = "'99.00"
$ws.Range("H38").Value = "'400.00"
$ws.Range("H39").Value = "'520.00"
$ws.Range("H40").Value = "'1530.00"
$ws.Range("H41").Value = "'51214.47"
$ws.Range("H42").Value = "'62375.35"
$ws.Range("H43").Value = "'758.88"
$ws.Range("H44").Value = "'26047.00"
$ws.Range("H45").Value = "'99430.00"
$ws.Range("H46").Value = "'508090.00"
$ws.Range("H47").Value = "'10639.35"
$ws.Range("H48").Value = "'7256.26"
$ws.Range("H49").Value = "'33930.00"
$ws.Range("H50").Value = "'44192.00"
$ws.Range("H51").Value = "'10565.52"
$ws.Range("H52").Value = "'186150.54"
$ws.Range("H53").Value = "'740.00"
$ws.Range("H54").Value = "'200.00"
$ws.Range("H55").Value = "'9680.00"
$ws.Range("H56").Value = "'184.50"
$ws.Range("H57").Value = "'300.23"
$ws.Range("H58").Value = "'3484.56"
$ws.Range("H59").Value = "'1764.25"
$ws.Range("H60").Value = "'8284.00"
$ws.Range("H61").Value = "'43513.02"
$ws.Range("H62").Value = "'130.05"
$ws.Range("H63").Value = "'9452.45"
$ws.Range("H64").Value = "'501.50"
$ws.Range("H65").Value = "'148142.89"
$ws.Range("H66").Value = "'8400.00"
$ws.Range("H67").Value = "'150694.62"
$ws.Range("H68").Value = "'7700.48"
$ws.Range("H69").Value = "'11770.84"
$ws.Range("H70").Value = "'35087.25"
$ws.Range("H71").Value = "'15880.00"
$ws.Range("H72").Value = "'4650.00"
$ws.Range("H73").Value = "'2750.00"
$ws.Range("H74").Value = "'4930.76"
$ws.Range("H75").Value = "'5050.00"
$ws.Range("H76").Value = "'8000.00"
$ws.Range("H77").Value = "'8750.00"
$ws.Range("H78").Value = "'2677.00"
$ws.Range("H79").Value = "'11783.60"
$ws.Range("H80").Value = "'9111.52"
$ws.Range("H81").Value = "'520.00"
$ws.Range("H82").Value = "'493.90"
$ws.Range("H83").Value = "'6115.55"
$ws.Range("H84").Value = "'600.00"
$ws.Range("H85").Value = "'1350.00"
$ws.Range("H86").Value = "'108000.00"
$ws.Range("H87").Value = "'480.00"
$ws.Range("H88").Value = "'300.00"
$ws.Range("H89").Value = "'4510.80"
$ws.Range("H90").Value = "'625.00"
$ws.Range("H91").Value = "'2560.00"
$ws.Range("H92").Value = "'18900.00"
$ws.Range("H93").Value = "'24000.00"
$ws.Range("H94").Value = "'7425.00"
$ws.Range("H95").Value = "'2600.00"
$ws.Range("H96").Value = "'378.00"
$ws.Range("H97").Value = "'26180.00"
$ws.Range("H98").Value = "'119.10"
$ws.Range("H99").Value = "'2900.00"
$ws.Range("H100").Value = "'21680.00"
$ws.Range("H101").Value = "'344.00"
$ws.Range("H102").Value = "'15108.00"
$ws.Range("H103").Value = "'2430.00"
$ws.Range("H104").Value = "'3350.00"
$ws.Range("H105").Value = "'558358.95"
$ws.Range("H106").Value = "'54317.09"
$ws.Range("H107").Value = "'44017080.99"
$ws.Range("H108").Value = "'2000.00"
$ws.Range("H109").Value = "'185.00"
$ws.Range("H110").Value = "'28703.62"
$ws.Range("H111").Value = "'17950.00"
$ws.Range("H112").Value = "'1920.00"
$ws.Range("H113").Value = "'10025.80"
$ws.Range("H114").Value = "'297.25"
$ws.Range("H115").Value = "'2636.28"
$ws.Range("H116").Value = "'3496.00"
$ws.Range("H117").Value = "'8550.00"
$ws.Range("H118").Value = "'24354.00"
$ws.Range("H119").Value = "'3240.00"
$ws.Range("H120").Value = "'6865.00"
$ws.Range("H121").Value = "'12160.00"
$ws.Range("H122").Value = "'34117.20"
$ws.Range("H123").Value = "'1540.00"
$ws.Range("H124").Value = "'23664.00"
$ws.Range("H125").Value = "'1180.00"
$ws.Range("H126").Value = "'5337.00"
$ws.Range("H127").Value = "'90.00"
$ws.Range("H128").Value = "'7993.00"
$ws.Range("H129").Value = "'4796.00"
$ws.Range("H130").Value = "'36368.77"
$ws.Range("H131").Value = "'10872.00"
$ws.Range("H132").Value = "'1757.05"
$ws.Range("H133").Value = "'674.90"
$ws.Range("H134").Value = "'16823.60"
$ws.Range("H135").Value = "'257.00"
$ws.Range("H136").Value = "'3300.00"
$ws.Range("H137").Value = "'184.10"
$ws.Range("H138").Value = "'66.41"
$ws.Range("H139").Value = "'1145.00"
$ws.Range("H140").Value = "'11200.00"
$ws.Range("H141").Value = "'150.00"
$ws.Range("H142").Value = "'612.00"
$ws.Range("H143").Value = "'810.00"
$ws.Range("H144").Value = "'1849.79"
$ws.Range("H145").Value = "'3750.00"
$ws.Range("H146").Value = "'439200.00"
$ws.Range("H147").Value = "'8000.00"
$ws.Range("H148").Value = "'497722.50"
$ws.Range("H149").Value = "'2820.00"
$ws.Range("H150").Value = "'3452.89"
$ws.Range("H151").Value = "'5400.00"
$ws.Range("H152").Value = "'3128.48"
$ws.Range("H153").Value = "'22508.00"
$ws.Range("H154").Value = "'692.00"
$ws.Range("H155").Value = "'3500.00"
$ws.Range("H156").Value = "'1100.00"
$ws.Range("H157").Value = "'9500.00"
$ws.Range("H158").Value = "'19000.00"
$ws.Range("H159").Value = "'2500.00"
$ws.Range("H160").Value = "'18700.00"
$ws.Range("H161").Value = "'9000.00"
$ws.Range("H162").Value = "'5851.34"
$ws.Range("H163").Value = "'4520.14"
$ws.Range("H164").Value = "'18135.10"
$ws.Range("H165").Value = "'1267.71"
$ws.Range("H166").Value = "'1719.00"
$ws.Range("H167").Value = "'14277.60"
$ws.Range("H168").Value = "'1093.50"
$ws.Range("H169").Value = "'3300000.00"
$ws.Range("H170").Value = "'356000.00"
$ws.Range("H171").Value = "'5080.00"
$ws.Range("H172").Value = "'1633.80"
$ws.Range("H173").Value = "'80628.00"
$ws.Range("H174").Value = "'10000.00"
$ws.Range("H175").Value = "'12000.00"
$ws.Range("H176").Value = "'14000.00"
$ws.Range("H177").Value = "'6000.00"
$ws.Range("H178").Value = "'48911.23"
$ws.Range("H179").Value = "'6500.00"
$ws.Range("H180").Value = "'3000.00"
$ws.Range("H181").Value = "'4999.12"
$ws.Range("H182").Value = "'4000.00"
$ws.Range("H183").Value = "'3000.00"
$ws.Range("H184").Value = "'6000.00"
$ws.Range("H185").Value = "'4500.00"
$ws.Range("H186").Value = "'3500.00"
$ws.Range("H187").Value = "'3000.00"
$ws.Range("H188").Value = "'9000.00"
$ws.Range("H189").Value = "'56465.50"
$ws.Range("H190").Value = "'10500.00"
$ws.Range("H191").Value = "'12000.00"
$ws.Range("H192").Value = "'9000.00"
$ws.Range("H193").Value = "'10100.00"
$ws.Range("H194").Value = "'21268.17"
$ws.Range("H195").Value = "'5000.00"
$ws.Range("H196").Value = "'45000.00"
$ws.Range("H197").Value = "'6000.00"
$ws.Range("H198").Value = "'4000.00"
$ws.Range("H199").Value = "'2500.00"
$ws.Range("H200").Value = "'75664.00"
$ws.Range("H201").Value = "'4000.00"
$ws.Range("H202").Value = "'21550.00"
$ws.Range("H203").Value = "'1200.00"
$ws.Range("H204").Value = "'300.00"
$ws.Range("H205").Value = "'5773.00"
$ws.Range("H206").Value = "'6151.00"
$ws.Range("H207").Value = "'6700.00"
$ws.Range("H208").Value = "'2637.60"
$ws.Range("H209").Value = "'636.00"
$ws.Range("H210").Value = "'4711.00"
$ws.Range("H211").Value = "'8645.00"
$ws.Range("H212").Value = "'12200.00"
$ws.Range("H213").Value = "'1913.50"
$ws.Range("H214").Value = "'1295.55"
$ws.Range("H215").Value = "'3770.00"
$ws.Range("H216").Value = "'921.00"
$ws.Range("H217").Value = "'9950.00"
$ws.Range("H218").Value = "'4929.50"
$ws.Range("H219").Value = "'1780.00"
$ws.Range("H220").Value = "'406.37"
$ws.Range("H221").Value = "'9808.00"
$ws.Range("H222").Value = "'5360.00"
$ws.Range("H223").Value = "'800.00"
$ws.Range("H224").Value = "'7410.00"
$ws.Range("H225").Value = "'31716.00"
$ws.Range("H226").Value = "'3874.05"
$ws.Range("H227").Value = "'804.96"
$ws.Range("H228").Value = "'1860.00"
$ws.Range("H229").Value = "'128500.00"
$ws.Range("H230").Value = "'500.00"
$ws.Range("H231").Value = "'448.26"
$ws.Range("H232").Value = "'2000.00"
$ws.Range("H233").Value = "'10710.00"
$ws.Range("H234").Value = "'2070.00"
$ws.Range("H235").Value = "'13600.00"
$ws.Range("H236").Value = "'4211.70"
$ws.Range("H237").Value = "'500.00"
$ws.Range("H238").Value = "'204271.48"
$ws.Range("H239").Value = "'8500.00"
$ws.Range("H240").Value = "'25000.00"
$ws.Range("H241").Value = "'50000.00"
$ws.Range("H242").Value = "'50000.00"
$ws.Range("H243").Value = "'25000.00"
$ws.Range("H244").Value = "'25000.00"
$ws.Range("H245").Value = "'50000.00"
$ws.Range("H246").Value = "'50000.00"
$ws.Range("H247").Value = "'50000.00"
$ws.Range("H248").Value = "'45600.00"
$ws.Range("H249").Value = "'1800.00"
$ws.Range("H250").Value = "'4000.00"
$ws.Range("H251").Value = "'2650250.00"
$ws.Range("H252").Value = "'274365.00"
$ws.Range("H253").Value = "'136500.00"
$ws.Range("H254").Value = "'135500.00"
$ws.Range("H255").Value = "'128000.00"
$ws.Range("H256").Value = "'128000.00"
$ws.Range("H257").Value = "'128000.00"
$ws.Range("H258").Value = "'128000.00"
$ws.Range("H259").Value = "'224000.00"
$ws.Range("H260").Value = "'224000.00"
$ws.Range("H261").Value = "'321500.00"
$ws.Range("H262").Value = "'128000.00"
$ws.Range("H263").Value = "'128000.00"
$ws.Range("H264").Value = "'128000.00"
$ws.Range("H265").Value = "'128000.00"
$ws.Range("H266").Value = "'128000.00"
$ws.Range("H267").Value = "'224000.00"
$ws.Range("H268").Value = "'320000.00"
$ws.Range("H269").Value = "'224000.00"
$ws.Range("H270").Value = "'128000.00"
$ws.Range("H271").Value = "'211000.00"
$ws.Range("H272").Value = "'128000.00"
$ws.Range("H273").Value = "'128000.00"
$ws.Range("H274").Value = "'143080.00"
$ws.Range("H275").Value = "'128000.00"
$ws.Range("H276").Value = "'2958038.42"
$ws.Range("H277").Value = "'244803.16"
$ws.Range("H278").Value = "'9156.00"
$ws.Range("H279").Value = "'12000.00"
$ws.Range("H280").Value = "'83000.00"
$ws.Range("H281").Value = "'30000.00"
$ws.Range("H282").Value = "'2783.00"
$ws.Range("H283").Value = "'9000.00"
$ws.Range("H284").Value = "'12500.00"
$ws.Range("H285").Value = "'4394.60"
$ws.Range("H286").Value = "'1200.00"
$ws.Range("H287").Value = "'7600.00"
$ws.Range("H288").Value = "'41400.00"
$ws.Range("H289").Value = "'900.00"
$ws.Range("H290").Value = "'100000.00"
$ws.Range("H291").Value = "'2800.00"
$ws.Range("H292").Value = "'7000.00"
$ws.Range("H293").Value = "'137986.00"
$ws.Range("H294").Value = "'5701.50"
$ws.Range("H295").Value = "'39800.00"

# Strip the quote-prefix marker left by the leading apostrophe so the
# cells end up with no explicit style (same as the original workbook).
$ws.Range("H2:H295").Style = "Normal"

